# Worked on SFD model, forces need work

$wb = $excel.ActiveWorkbook

# --- Vehicle Sections sheet: add new Fins-row data in columns G:N ---
$vehicleSections = $wb.Worksheets.Item("Vehicle Sections")
$vehicleSections.Range("G11").Value = 1
$vehicleSections.Range("H11").Value = 1
$vehicleSections.Range("I11").Value = 1
$vehicleSections.Range("J11").Value = 1
$vehicleSections.Range("K11").Value = 1
$vehicleSections.Range("L11").Value = 1
$vehicleSections.Range("M11").Value = 1
$vehicleSections.Range("N11").Value = 3

# Move the selection on "Vehicle Sections" off of D14 and onto C5
$vehicleSections.Range("C5").Select()

# --- New "Aerodynamic Properties" sheet, inserted right after "Vehicle Sections" ---
$newSheet = $wb.Worksheets.Add($null, $vehicleSections)
$newSheet.Name = "Aerodynamic Properties"

# Write the row-label strings first so they land earlier in the shared
# string table (matches "Off the rail" / "Max Q" preceding the column
# headers).
$newSheet.Range("A2").Value = "Off the rail"
$newSheet.Range("A3").Value = "Max Q"

$newSheet.Range("A1").Value = "Name"
$newSheet.Range("B1").Value = "Velocity (m/s)"
$newSheet.Range("C1").Value = "Acceleration (m/s^2)"
$newSheet.Range("D1").Value = "Mach"
$newSheet.Range("E1").Value = "Cd"
$newSheet.Range("F1").Value = "Thrust (lbf)"
$newSheet.Range("G1").Value = "Max wind gust (mph)"
$newSheet.Range("H1").Value = "AoA"
$newSheet.Range("I1").Value = "Air Density"

$newSheet.Range("B2").Value = 29.7
$newSheet.Range("C2").Value = 1
$newSheet.Range("D2").Value = 1
$newSheet.Range("E2").Value = 1
$newSheet.Range("F2").Value = 1
$newSheet.Range("G2").Value = 1
$newSheet.Range("H2").Value = 1
$newSheet.Range("I2").Value = 1

$newSheet.Range("B3").Value = 1
$newSheet.Range("C3").Value = 1
$newSheet.Range("D3").Value = 1
$newSheet.Range("E3").Value = 1
$newSheet.Range("F3").Value = 1
$newSheet.Range("G3").Value = 1
$newSheet.Range("H3").Value = 1
$newSheet.Range("I3").Value = 1

$newSheet.Range("C3").Select()

# Activate the new sheet so it becomes the active tab (index 1)
$newSheet.Activate()
